# Generate Report for Handback
# Adds a new handback row (file 1005d520-1be8-469c-8bf4-db190a8f2e97) to the
# Overview, zh-cn and de-de tables/sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table "Overview") -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.md"
$wsOverview.Range("B3").Value = "e2e\1005d520-1be8-469c-8bf4-db190a8f2e97.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G3").Value = "2016-12-15 03:49:19"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/1005d520-1be8-469c-8bf4-db190a8f2e97.md", "", "", "e2e\1005d520-1be8-469c-8bf4-db190a8f2e97.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table "zh-cn") -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.fcd5dd8dd2d9fb862dc4b854ba4d7d5c77ae1c49.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "2016-12-15 03:49:05"
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("I3").Style = "Normal"
$wsZhCn.Range("J3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.md"
$wsZhCn.Range("K3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.fcd5dd8dd2d9fb862dc4b854ba4d7d5c77ae1c49.zh-cn.xlf"
$wsZhCn.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = "2016-12-15 03:50:00"
$wsZhCn.Range("M3").Value = "'"
$wsZhCn.Range("M3").Style = "Normal"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("N3").Style = "Normal"
$wsZhCn.Range("O3").Value = "'True"
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("P3").Style = "Normal"
$wsZhCn.Range("Q3").Value = "'False"
$wsZhCn.Range("R3").Value = "'"
$wsZhCn.Range("R3").Style = "Normal"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/1005d520-1be8-469c-8bf4-db190a8f2e97.md", "", "", "1005d520-1be8-469c-8bf4-db190a8f2e97.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/1005d520-1be8-469c-8bf4-db190a8f2e97.md", "", "", "1005d520-1be8-469c-8bf4-db190a8f2e97.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table "de-de") -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.fcd5dd8dd2d9fb862dc4b854ba4d7d5c77ae1c49.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "2016-12-15 03:49:19"
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("I3").Style = "Normal"
$wsDeDe.Range("J3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.md"
$wsDeDe.Range("K3").Value = "1005d520-1be8-469c-8bf4-db190a8f2e97.fcd5dd8dd2d9fb862dc4b854ba4d7d5c77ae1c49.de-de.xlf"
$wsDeDe.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = "2016-12-15 03:50:19"
$wsDeDe.Range("M3").Value = "'"
$wsDeDe.Range("M3").Style = "Normal"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("N3").Style = "Normal"
$wsDeDe.Range("O3").Value = "'True"
$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Range("P3").Style = "Normal"
$wsDeDe.Range("Q3").Value = "'False"
$wsDeDe.Range("R3").Value = "'"
$wsDeDe.Range("R3").Style = "Normal"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/1005d520-1be8-469c-8bf4-db190a8f2e97.md", "", "", "1005d520-1be8-469c-8bf4-db190a8f2e97.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/1005d520-1be8-469c-8bf4-db190a8f2e97.md", "", "", "1005d520-1be8-469c-8bf4-db190a8f2e97.md") | Out-Null

$wb.Save()
